$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.782436333333333
$ws.Cells.Item(2, 8).Value = 5.347308999999999
$ws.Cells.Item(2, 9).Value = 0.00914036392049929
$ws.Cells.Item(2, 10).Value = 0.009140363920499292
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.4702473333333333
$ws.Cells.Item(2, 14).Value = 1.410742
$ws.Cells.Item(2, 15).Value = 0.00903492226842282
$ws.Cells.Item(2, 16).Value = 0.00903492226842282
$ws.Cells.Item(2, 17).Value = 0.8381859325864442
$ws.Cells.Item(2, 18).Value = 7.543673393277999
$ws.Cells.Item(2, 19).Value = 0.00008258247752680754
$ws.Cells.Item(2, 20).Value = 0.00008258247752680756
$ws.Cells.Item(3, 7).Value = 1.782436333333333
$ws.Cells.Item(3, 8).Value = 5.347308999999999
$ws.Cells.Item(3, 9).Value = 0.00914036392049929
$ws.Cells.Item(3, 10).Value = 0.009140363920499292
$ws.Cells.Item(3, 14).Value = 0.9584440000000001
$ws.Cells.Item(3, 15).Value = 0.006138235792679485
$ws.Cells.Item(3, 16).Value = 0.006138235792679485
$ws.Cells.Item(3, 17).Value = 0.569455136355111
$ws.Cells.Item(3, 18).Value = 5.125096227196
$ws.Cells.Item(3, 19).Value = 0.00005610570897492492
$ws.Cells.Item(3, 20).Value = 0.00005610570897492494
$ws.Cells.Item(4, 7).Value = 1.782436333333333
$ws.Cells.Item(4, 8).Value = 5.347308999999999
$ws.Cells.Item(4, 9).Value = 0.00914036392049929
$ws.Cells.Item(4, 10).Value = 0.009140363920499292
$ws.Cells.Item(4, 13).Value = 1.047307
$ws.Cells.Item(4, 14).Value = 3.141921
$ws.Cells.Item(4, 15).Value = 0.02012204358311108
$ws.Cells.Item(4, 16).Value = 0.02012204358311108
$ws.Cells.Item(4, 17).Value = 1.866758048954333
$ws.Cells.Item(4, 18).Value = 16.800822440589
$ws.Cells.Item(4, 19).Value = 0.0001839228011737828
$ws.Cells.Item(4, 20).Value = 0.0001839228011737828
$ws.Cells.Item(5, 7).Value = 1.782436333333333
$ws.Cells.Item(5, 8).Value = 5.347308999999999
$ws.Cells.Item(5, 9).Value = 0.00914036392049929
$ws.Cells.Item(5, 10).Value = 0.009140363920499292
$ws.Cells.Item(5, 13).Value = 50.21070966666667
$ws.Cells.Item(5, 14).Value = 150.632129
$ws.Cells.Item(5, 15).Value = 0.9647047983557866
$ws.Cells.Item(5, 16).Value = 0.9647047983557866
$ws.Cells.Item(5, 17).Value = 89.49739323231789
$ws.Cells.Item(5, 18).Value = 805.476539090861
$ws.Cells.Item(5, 19).Value = 0.008817752932823775
$ws.Cells.Item(5, 20).Value = 0.008817752932823776
$ws.Cells.Item(6, 8).Value = 564.692825
$ws.Cells.Item(6, 9).Value = 0.965251479537618
$ws.Cells.Item(6, 10).Value = 0.965251479537618
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.4702473333333333
$ws.Cells.Item(6, 14).Value = 1.410742
$ws.Cells.Item(6, 15).Value = 0.00903492226842282
$ws.Cells.Item(6, 16).Value = 0.00903492226842282
$ws.Cells.Item(6, 17).Value = 88.51509836957221
$ws.Cells.Item(6, 18).Value = 796.63588532615
$ws.Cells.Item(6, 19).Value = 0.0087209720871025
$ws.Cells.Item(6, 20).Value = 0.0087209720871025
$ws.Cells.Item(7, 8).Value = 564.692825
$ws.Cells.Item(7, 9).Value = 0.965251479537618
$ws.Cells.Item(7, 10).Value = 0.965251479537618
$ws.Cells.Item(7, 14).Value = 0.9584440000000001
$ws.Cells.Item(7, 15).Value = 0.006138235792679485
$ws.Cells.Item(7, 16).Value = 0.006138235792679485
$ws.Cells.Item(7, 18).Value = 541.2264499643001
$ws.Cells.Item(7, 19).Value = 0.005924941180634637
$ws.Cells.Item(7, 20).Value = 0.005924941180634637
$ws.Cells.Item(8, 8).Value = 564.692825
$ws.Cells.Item(8, 9).Value = 0.965251479537618
$ws.Cells.Item(8, 10).Value = 0.965251479537618
$ws.Cells.Item(8, 13).Value = 1.047307
$ws.Cells.Item(8, 14).Value = 3.141921
$ws.Cells.Item(8, 15).Value = 0.02012204358311108
$ws.Cells.Item(8, 16).Value = 0.02012204358311108
$ws.Cells.Item(8, 17).Value = 197.1355828240917
$ws.Cells.Item(8, 18).Value = 1774.220245416825
$ws.Cells.Item(8, 19).Value = 0.01942283233991841
$ws.Cells.Item(8, 20).Value = 0.01942283233991841
$ws.Cells.Item(9, 8).Value = 564.692825
$ws.Cells.Item(9, 9).Value = 0.965251479537618
$ws.Cells.Item(9, 10).Value = 0.965251479537618
$ws.Cells.Item(9, 13).Value = 50.21070966666667
$ws.Cells.Item(9, 14).Value = 150.632129
$ws.Cells.Item(9, 15).Value = 0.9647047983557866
$ws.Cells.Item(9, 16).Value = 0.9647047983557866
$ws.Cells.Item(9, 17).Value = 9451.20916230827
$ws.Cells.Item(9, 18).Value = 85060.88246077443
$ws.Cells.Item(9, 19).Value = 0.9311827339299624
$ws.Cells.Item(9, 20).Value = 0.9311827339299624
$ws.Cells.Item(10, 7).Value = 4.870778333333333
$ws.Cells.Item(10, 8).Value = 14.612335
$ws.Cells.Item(10, 9).Value = 0.02497743437460768
$ws.Cells.Item(10, 10).Value = 0.02497743437460768
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4702473333333333
$ws.Cells.Item(10, 14).Value = 1.410742
$ws.Cells.Item(10, 15).Value = 0.00903492226842282
$ws.Cells.Item(10, 16).Value = 0.00903492226842282
$ws.Cells.Item(10, 17).Value = 2.290470522507777
$ws.Cells.Item(10, 18).Value = 20.61423470256999
$ws.Cells.Item(10, 19).Value = 0.0002256691780392125
$ws.Cells.Item(10, 20).Value = 0.0002256691780392125
$ws.Cells.Item(11, 7).Value = 4.870778333333333
$ws.Cells.Item(11, 8).Value = 14.612335
$ws.Cells.Item(11, 9).Value = 0.02497743437460768
$ws.Cells.Item(11, 10).Value = 0.02497743437460768
$ws.Cells.Item(11, 14).Value = 0.9584440000000001
$ws.Cells.Item(11, 15).Value = 0.006138235792679485
$ws.Cells.Item(11, 16).Value = 0.006138235792679485
$ws.Cells.Item(11, 17).Value = 1.556122756304444
$ws.Cells.Item(11, 18).Value = 14.00510480674
$ws.Cells.Item(11, 19).Value = 0.0001533173816875198
$ws.Cells.Item(11, 20).Value = 0.0001533173816875198
$ws.Cells.Item(12, 7).Value = 4.870778333333333
$ws.Cells.Item(12, 8).Value = 14.612335
$ws.Cells.Item(12, 9).Value = 0.02497743437460768
$ws.Cells.Item(12, 10).Value = 0.02497743437460768
$ws.Cells.Item(12, 13).Value = 1.047307
$ws.Cells.Item(12, 14).Value = 3.141921
$ws.Cells.Item(12, 15).Value = 0.02012204358311108
$ws.Cells.Item(12, 16).Value = 0.02012204358311108
$ws.Cells.Item(12, 17).Value = 5.101200243948333
$ws.Cells.Item(12, 18).Value = 45.91080219553499
$ws.Cells.Item(12, 19).Value = 0.0005025970230801526
$ws.Cells.Item(12, 20).Value = 0.0005025970230801526
$ws.Cells.Item(13, 7).Value = 4.870778333333333
$ws.Cells.Item(13, 8).Value = 14.612335
$ws.Cells.Item(13, 9).Value = 0.02497743437460768
$ws.Cells.Item(13, 10).Value = 0.02497743437460768
$ws.Cells.Item(13, 13).Value = 50.21070966666667
$ws.Cells.Item(13, 14).Value = 150.632129
$ws.Cells.Item(13, 15).Value = 0.9647047983557866
$ws.Cells.Item(13, 16).Value = 0.9647047983557866
$ws.Cells.Item(13, 17).Value = 244.5652367456906
$ws.Cells.Item(13, 18).Value = 2201.087130711215
$ws.Cells.Item(13, 19).Value = 0.02409585079180079
$ws.Cells.Item(13, 20).Value = 0.02409585079180079
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.1229953333333333
$ws.Cells.Item(14, 8).Value = 0.368986
$ws.Cells.Item(14, 9).Value = 0.0006307221672750447
$ws.Cells.Item(14, 10).Value = 0.0006307221672750447
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.4702473333333333
$ws.Cells.Item(14, 14).Value = 1.410742
$ws.Cells.Item(14, 15).Value = 0.00903492226842282
$ws.Cells.Item(14, 16).Value = 0.00903492226842282
$ws.Cells.Item(14, 17).Value = 0.05783822751244445
$ws.Cells.Item(14, 18).Value = 0.520544047612
$ws.Cells.Item(14, 19).Value = 0.000005698525754301204
$ws.Cells.Item(14, 20).Value = 0.000005698525754301204
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.1229953333333333
$ws.Cells.Item(15, 8).Value = 0.368986
$ws.Cells.Item(15, 9).Value = 0.0006307221672750447
$ws.Cells.Item(15, 10).Value = 0.0006307221672750447
$ws.Cells.Item(15, 14).Value = 0.9584440000000001
$ws.Cells.Item(15, 15).Value = 0.006138235792679485
$ws.Cells.Item(15, 16).Value = 0.006138235792679485
$ws.Cells.Item(15, 17).Value = 0.03929471308711112
$ws.Cells.Item(15, 18).Value = 0.3536524177840001
$ws.Cells.Item(15, 19).Value = 0.000003871521382404057
$ws.Cells.Item(15, 20).Value = 0.000003871521382404057
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.1229953333333333
$ws.Cells.Item(16, 8).Value = 0.368986
$ws.Cells.Item(16, 9).Value = 0.0006307221672750447
$ws.Cells.Item(16, 10).Value = 0.0006307221672750447
$ws.Cells.Item(16, 13).Value = 1.047307
$ws.Cells.Item(16, 14).Value = 3.141921
$ws.Cells.Item(16, 15).Value = 0.02012204358311108
$ws.Cells.Item(16, 16).Value = 0.02012204358311108
$ws.Cells.Item(16, 17).Value = 0.1288138735673333
$ws.Cells.Item(16, 18).Value = 1.159324862106
$ws.Cells.Item(16, 19).Value = 0.00001269141893874273
$ws.Cells.Item(16, 20).Value = 0.00001269141893874273
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.1229953333333333
$ws.Cells.Item(17, 8).Value = 0.368986
$ws.Cells.Item(17, 9).Value = 0.0006307221672750447
$ws.Cells.Item(17, 10).Value = 0.0006307221672750447
$ws.Cells.Item(17, 13).Value = 50.21070966666667
$ws.Cells.Item(17, 14).Value = 150.632129
$ws.Cells.Item(17, 15).Value = 0.9647047983557866
$ws.Cells.Item(17, 16).Value = 0.9647047983557866
$ws.Cells.Item(17, 17).Value = 6.17568297235489
$ws.Cells.Item(17, 18).Value = 55.58114675119401
$ws.Cells.Item(17, 19).Value = 0.0006084607011995968
$ws.Cells.Item(17, 20).Value = 0.0006084607011995968
